$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark cells whose text values look numeric as Text format,
# so COM does not silently coerce them into numbers (matches the
# original workbook, where these columns are stored as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '62.493.07'
$ws.Range("E2").Value = '  -0.95%  '

$ws.Range("D3").Value = '3.010.15'
$ws.Range("E3").Value = '  -1.10%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '585.57'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").Value = '147.10'
$ws.Range("E6").Value = '  -3.02%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  -2.27%  '

$ws.Range("D9").Value = '3.008.65'
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("E10").Value = '  -3.16%  '

$ws.Range("D11").Value = '5.81'
$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("E12").Value = '  +3.28%  '

$ws.Range("D13").Value = '0.0000230'
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").Value = '34.65'
$ws.Range("E14").Value = '  -4.72%  '

$ws.Range("E15").Value = '  +2.37%  '

$ws.Range("D16").Value = '3.513.50'
$ws.Range("E16").Value = '  -0.81%  '

$ws.Range("D17").Value = '7.09'
$ws.Range("E17").Value = '  -0.93%  '

$ws.Range("D18").Value = '62.473.32'
$ws.Range("E18").Value = '  -0.94%  '

$ws.Range("D19").Value = '3.012.60'
$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("D20").Value = '459.52'
$ws.Range("E20").Value = '  -3.93%  '

$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  -1.80%  '

$ws.Range("D22").Value = '0.690'
$ws.Range("E22").Value = '  -2.25%  '

$ws.Range("D23").Value = '7.43'
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("D24").Value = '81.76'
$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -8.23%  '

$ws.Range("D26").Value = '12.22'
$ws.Range("E26").Value = '  -3.80%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").Value = '9.89'
$ws.Range("E28").Value = '  -7.59%  '

$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").Value = '2.65'
$ws.Range("E30").Value = '  -0.73%  '

$ws.Range("D31").Value = '7.00'
$ws.Range("E31").Value = '  -5.02%  '

$ws.Range("D32").Value = '2.10'
$ws.Range("E32").Value = '  -4.54%  '

$ws.Range("D33").Value = '27.87'
$ws.Range("E33").Value = '  +0.81%  '

$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  -1.37%  '

$ws.Range("D35").Value = '0.0₃0814'
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  -3.02%  '

$ws.Range("D37").Value = '5.77'
$ws.Range("E37").Value = '  -2.45%  '

$ws.Range("E38").Value = '  -4.53%  '

$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").Value = '9.19'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '50.39'
$ws.Range("E40").Value = '  -0.13%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.123'
$ws.Range("E41").Value = '  +8.77%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '2.91'
$ws.Range("E42").Value = '  -10.60%  '

$ws.Range("D43").Value = '391.56'
$ws.Range("E43").Value = '  -10.11%  '

$ws.Range("D44").Value = '0.0359'
$ws.Range("E44").Value = '  -0.78%  '

$ws.Range("D45").Value = '0.268'
$ws.Range("E45").Value = '  -6.38%  '

$ws.Range("D46").Value = '2.742.63'
$ws.Range("E46").Value = '  -2.85%  '

$ws.Range("D47").Value = '37.41'
$ws.Range("E47").Value = '  -2.48%  '

$ws.Range("D48").Value = '129.36'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("E49").Value = '  +0.09%  '

$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("E51").Value = '  -0.75%  '
